# Update "Inscritos" / "Pagos" / "Inscrições homologadas" counts in the
# Table1 range (sheet "Inscricoes") to reflect newly registered applicants.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 85

$ws.Range("E4").Value = 39

$ws.Range("E5").Value = 118

$ws.Range("E6").Value = 39
$ws.Range("F6").Value = 26
$ws.Range("H6").Value = 26

$ws.Range("E10").Value = 393
$ws.Range("F10").Value = 188
$ws.Range("H10").Value = 188

$ws.Range("E11").Value = 260
$ws.Range("F11").Value = 140
$ws.Range("H11").Value = 140

$ws.Range("E12").Value = 378
$ws.Range("F12").Value = 211
$ws.Range("H12").Value = 211

$ws.Range("E13").Value = 103

$ws.Range("E16").Value = 168

$ws.Range("E17").Value = 74

$ws.Range("E18").Value = 48
$ws.Range("F18").Value = 21
$ws.Range("H18").Value = 21

$ws.Range("E20").Value = 76

$ws.Range("E21").Value = 123

$ws.Range("E22").Value = 144

$ws.Range("E23").Value = 168

$ws.Range("E24").Value = 168
$ws.Range("F24").Value = 86
$ws.Range("H24").Value = 86

$ws.Range("E25").Value = 198
$ws.Range("F25").Value = 90
$ws.Range("H25").Value = 90

$ws.Range("E27").Value = 261
$ws.Range("F27").Value = 124
$ws.Range("H27").Value = 124

$ws.Range("E30").Value = 169

$ws.Range("E31").Value = 65

$ws.Range("E32").Value = 157
$ws.Range("F32").Value = 89
$ws.Range("H32").Value = 89

$ws.Range("E33").Value = 240

$ws.Range("E34").Value = 172
$ws.Range("F34").Value = 103
$ws.Range("H34").Value = 103

$ws.Range("E35").Value = 114

$ws.Range("E37").Value = 127

$ws.Range("E39").Value = 159

$ws.Range("E40").Value = 215

$ws.Range("E41").Value = 316
$ws.Range("F41").Value = 138
$ws.Range("H41").Value = 138

$ws.Range("E42").Value = 286

$ws.Range("E43").Value = 96
$ws.Range("F43").Value = 49
$ws.Range("H43").Value = 49

$ws.Range("E44").Value = 252

$ws.Range("E45").Value = 111
$ws.Range("F45").Value = 51
$ws.Range("H45").Value = 51

$ws.Range("E46").Value = 246

$ws.Range("E47").Value = 355
$ws.Range("F47").Value = 174
$ws.Range("H47").Value = 174

$ws.Range("E48").Value = 167

$ws.Range("E49").Value = 241
$ws.Range("F49").Value = 102
$ws.Range("H49").Value = 102

$ws.Range("F50").Value = 84
$ws.Range("H50").Value = 84

$ws.Range("E51").Value = 193
